$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM")

# --- Update VALIDATIONS (column J) for the password-change test rows. ---
# The API responses used to echo back the HTTP status code as `errorcode`
# (e.g. errorcode=422); the service now always reports `errorcode=0`.
$ws.Range("J38").Value = "status=422||errorcode=0||reason=New password should not match current password"
$ws.Range("J39").Value = "status=422||errorcode=0||reason=New password should not match previous 4 passwords"
$ws.Range("J41").Value = "status=400||errorcode=0||reason=Update request body is missing required parameters"
$ws.Range("J42").Value = "status=400||errorcode=0||reason=Update request body is missing required parameters"
$ws.Range("J43").Value = "status=422||errorcode=0"
$ws.Range("J44").Value = "status=422||errorcode=0||reason=Password should be at least 8 characters long||reason=Password should contain at least one alphabet character, either upper or lower case"
$ws.Range("J45").Value = "status=422||errorcode=0||reason=Password should have at least 1 numeric character"
$ws.Range("J46").Value = "status=422||errorcode=0||reason=Password should be at least 8 characters long||reason=Password should have at least 1 numeric character"

# --- Move the sheet's saved view back to the top of the sheet and select H3. ---
$ws.Activate() | Out-Null
$ws.Range("H3").Select() | Out-Null
